# Added new data from another group for HVO60
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate Sheet1 (keeps all formatting/merged cells/styles identical) and
# place the copy right after Sheet1, then rename it to Sheet2.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"

# New HVO60 measurement data for the second group, rows 22-28 (columns B:G).
# Column A (the "CA / reading" index 14..20) stays identical to Sheet1.
$ws2.Range("B22").Value = 0
$ws2.Range("C22").Value = 0
$ws2.Range("D22").Value = 335
$ws2.Range("E22").Value = 4.5
$ws2.Range("F22").Value = 14.42
$ws2.Range("G22").Value = "-"

$ws2.Range("B23").Value = 0
$ws2.Range("C23").Value = 0
$ws2.Range("D23").Value = 360
$ws2.Range("E23").Value = 4.5
$ws2.Range("F23").Value = 14.45
$ws2.Range("G23").Value = "-"

$ws2.Range("B24").Value = 0
$ws2.Range("C24").Value = 0
$ws2.Range("D24").Value = 385
$ws2.Range("E24").Value = 45
$ws2.Range("F24").Value = 14.44
$ws2.Range("G24").Value = "-"

$ws2.Range("B25").Value = 0
$ws2.Range("C25").Value = 0
$ws2.Range("D25").Value = 424
$ws2.Range("E25").Value = 4.5
$ws2.Range("F25").Value = 14.4
$ws2.Range("G25").Value = "-"

$ws2.Range("B26").Value = 0
$ws2.Range("C26").Value = 0
$ws2.Range("D26").Value = 477
$ws2.Range("E26").Value = 4.5999999999999996
$ws2.Range("F26").Value = 14.2
$ws2.Range("G26").Value = "-"

$ws2.Range("B27").Value = 0
$ws2.Range("C27").Value = 0
$ws2.Range("D27").Value = 0
$ws2.Range("E27").Value = 4.4000000000000004
$ws2.Range("F27").Value = 14.45
$ws2.Range("G27").Value = "-"

$ws2.Range("B28").Value = 0
$ws2.Range("C28").Value = 0
$ws2.Range("D28").Value = 0
$ws2.Range("E28").Value = 4.5
$ws2.Range("F28").Value = 14.38
$ws2.Range("G28").Value = "-"

# Sheet1 ends up deselecting its tab and keeps its scroll position (A7),
# but the live selection moved to I27.
$ws1.Range("I27").Select() | Out-Null

# Sheet2 becomes the active sheet/tab; its view is scrolled down to the
# second table and the live selection sits on J28 (first empty cell below
# the new data).
$ws2.Activate() | Out-Null
$ws2.Range("J28").Select() | Out-Null
